$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.143.89'
$ws.Range('E2').Value = '  +1.63%  '
$ws.Range('D3').Value = '1.998.72'
$ws.Range('E3').Value = '  +2.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.08%  '
$ws.Range('E6').Value = '  +2.76%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.25'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.87%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +2.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0801'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.50%  '
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.91'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.79%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.66'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.74%  '
$ws.Range('D14').Value = '2.293.22'
$ws.Range('E14').Value = '  +2.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.846'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.46'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.36%  '
$ws.Range('D17').Value = '2.001.36'
$ws.Range('E17').Value = '  +2.32%  '
$ws.Range('D18').Value = '37.058.79'
$ws.Range('E18').Value = '  +1.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '70.26'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.00%  '
$ws.Range('E20').Value = '  +1.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.18'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '230.38'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.48'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.21%  '
$ws.Range('E25').Value = '  +0.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.39'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.96%  '
$ws.Range('E27').Value = '  +5.87%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '163.27'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.97%  '
$ws.Range('E29').Value = '  +1.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.34'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +15.07%  '
$ws.Range('E31').Value = '  +1.85%  '
$ws.Range('E32').Value = '  +2.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0657'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +7.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.52'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.37%  '
$ws.Range('E35').Value = '  +5.65%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('E37').Value = '  +2.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.27'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.63%  '
$ws.Range('E39').Value = '  +3.97%  '
$ws.Range('E40').Value = '  +0.44%  '
$ws.Range('E41').Value = '  +0.54%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0214'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.52%  '
$ws.Range('E43').Value = '  +2.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.64'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '90.93'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.06%  '
$ws.Range('D46').Value = '1.379.15'
$ws.Range('E46').Value = '  +1.27%  '
$ws.Range('E47').Value = '  +2.55%  '
$ws.Range('E48').Value = '  +1.52%  '
$ws.Range('E49').Value = '  +0.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.02'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +15.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '46.27'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.57%  '
